# Updated cryptos list - apply new coin rankings, prices and volume values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.699.43'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.65%  '

$ws.Range("B3").Value = 'Ethereum'
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.949.41'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.86%  '

$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.000'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("B5").Value = 'BNB'
$ws.Range("C5").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '247.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.08%  '

$ws.Range("B6").Value = 'USDC'
$ws.Range("C6").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("B7").Value = 'XRP'
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4835'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.01%  '

$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2941'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.71%  '

$ws.Range("B9").Value = 'Dogecoin'
$ws.Range("C9").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06830'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.88%  '

$ws.Range("B10").Value = 'Litecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '112.47'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.64%  '

$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.43'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.94%  '

$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.953.53'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.09%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07668'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.45%  '

$ws.Range("B14").Value = 'Polkadot'
$ws.Range("C14").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.525'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.89%  '

$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6916'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.50%  '

$ws.Range("B16").Value = 'BitcoinCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '296.39'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +7.79%  '

$ws.Range("B17").Value = 'WrappedBTC'
$ws.Range("C17").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.721.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.75%  '

$ws.Range("B18").Value = 'Avalanche'
$ws.Range("C18").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.30'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.50%  '

$ws.Range("B19").Value = 'Uniswap'
$ws.Range("C19").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.682'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +3.62%  '

$ws.Range("B20").Value = 'ShibaInu'
$ws.Range("C20").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007720'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.50%  '

$ws.Range("B21").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C21").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.201.54'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.83%  '

$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.02%  '

$ws.Range("B23").Value = 'BinanceUSD'
$ws.Range("C23").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.04%  '

$ws.Range("B24").Value = 'Chainlink'
$ws.Range("C24").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.599'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.27%  '

$ws.Range("B25").Value = 'Cosmos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.864'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.51%  '

$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.47'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.94%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.85%  '

$ws.Range("B28").Value = 'LidoDAOToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.192'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.33%  '

$ws.Range("B29").Value = 'Stellar'
$ws.Range("C29").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1091'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.81%  '

$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.439'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.97%  '

$ws.Range("B31").Value = 'Filecoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.759'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +17.48%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.463'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +7.77%  '

$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05097'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.06%  '

$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7775'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.72%  '

$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.161'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.58%  '

$ws.Range("B36").Value = 'VeChain'
$ws.Range("C36").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02083'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.65%  '

$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.731'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.07%  '

$ws.Range("B38").Value = 'MXToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.698'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.09%  '

$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.057'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.81%  '

$ws.Range("B40").Value = 'Quant'
$ws.Range("C40").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '111.16'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.31%  '

$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.4478'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.28%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8746'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.99%  '

$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.923'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.29%  '

$ws.Range("B44").Value = 'Aave'
$ws.Range("C44").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '70.25'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.51%  '

$ws.Range("B45").Value = 'PaxDollar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.003'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.25%  '

$ws.Range("B46").Value = 'Aptos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.383'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.50%  '

$ws.Range("B47").Value = 'EnergySwap'
$ws.Range("C47").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.376'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.99%  '

$ws.Range("B48").Value = 'BitcoinSV'
$ws.Range("C48").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '48.70'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.10%  '

$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1256'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.08%  '

$ws.Range("B50").Value = 'Elrond'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '35.71'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +2.94%  '

$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.2541'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.72%  '

Write-Host "Cryptos list updated"